$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.232.34"
$ws.Range("E2").Value = "  -4.71%  "

$ws.Range("D3").Value = "2.232.88"
$ws.Range("E3").Value = "  -5.74%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.587"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.69%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.561"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0827"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.34%  "

$ws.Range("E14").Value = "  -1.61%  "

$ws.Range("D15").Value = "2.568.81"
$ws.Range("E15").Value = "  -5.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.861"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -12.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.21%  "

$ws.Range("D18").Value = "2.235.38"
$ws.Range("E18").Value = "  -5.36%  "

$ws.Range("D19").Value = "43.017.33"
$ws.Range("E19").Value = "  -5.16%  "

$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "0.0₃0964"
$ws.Range("E21").Value = "  -9.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "238.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.27%  "

$ws.Range("E26").Value = "  -9.66%  "

$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.11%  "

$ws.Range("E30").Value = "  -2.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0872"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.35%  "

$ws.Range("E36").Value = "  -5.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.122"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.68%  "

$ws.Range("E40").Value = "  -5.71%  "

$ws.Range("E41").Value = "  -11.34%  "

$ws.Range("E42").Value = "  -5.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0323"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.77%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "1.800.21"
$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.206"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.22%  "
